$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write new cells in the exact order the strings were first introduced ---
# (preserves shared-string table ordering: pres, P03, P04, P05, esito, P06, P07,
# hospital_sec_diagnosis)
$ws.Range("H1").Value = "pres"
$ws.Range("A5").Value = "P03"
$ws.Range("A6").Value = "P04"
$ws.Range("A7").Value = "P05"
$ws.Range("G1").Value = "esito"
$ws.Range("A9").Value = "P06"
$ws.Range("A10").Value = "P07"
$ws.Range("F4").Value = "hospital_sec_diagnosis"

$ws.Range("G2").Value = 2
$ws.Range("H3").Value = 1

$ws.Range("A4").Value = "P01"
$ws.Range("B4").Value = 200
$ws.Range("E4").Value = "sdo"
$ws.Range("H4").Value = 1

$ws.Range("B5").Value = 250
$ws.Range("E5").Value = "ps"
$ws.Range("F5").Value = "emergency_room_diagnosis"
$ws.Range("G5").Value = 1

$ws.Range("B6").Value = 200
$ws.Range("E6").Value = "sdo"
$ws.Range("F6").Value = "hospital_sec_diagnosis"
$ws.Range("H6").Value = 0

$ws.Range("B7").Value = 300
$ws.Range("E7").Value = "sdo"
$ws.Range("F7").Value = "hospital_main_diagnosis"
$ws.Range("H7").Value = 0

$ws.Range("A8").Value = "P04"
$ws.Range("B8").Value = 205
$ws.Range("E8").Value = "sdo"
$ws.Range("F8").Value = "hospital_sec_diagnosis"

$ws.Range("B9").Value = 600
$ws.Range("E9").Value = "sdo"
$ws.Range("F9").Value = "hospital_main_diagnosis"
$ws.Range("H9").Value = 0

$ws.Range("B10").Value = 250
$ws.Range("E10").Value = "sdo"
$ws.Range("F10").Value = "hospital_main_diagnosis"

# --- Page setup (adds pageSetup element: paperSize=9 / portrait) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Move the active selection to match the post-edit cursor position ---
$ws.Range("F8").Select()
